# Data refresh for the cryptos worksheet: updates Price (D) / Volume(1h) (E)
# text values, and swaps the RenderToken/Dai rows (27 <-> 28) including their
# Coin name, Link and Price/Volume cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'62.320.59"
$ws.Range('E2').Value = '  -3.03%  '

# Row 3
$ws.Range('D3').Value = "'3.000.00"
$ws.Range('E3').Value = '  -3.90%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').Value = "'581.14"
$ws.Range('E5').Value = '  -2.26%  '

# Row 6
$ws.Range('D6').Value = "'146.27"
$ws.Range('E6').Value = '  -7.46%  '

# Row 7
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('D8').Value = "'0.523"
$ws.Range('E8').Value = '  -3.52%  '

# Row 9
$ws.Range('D9').Value = "'2.999.61"
$ws.Range('E9').Value = '  -3.89%  '

# Row 10
$ws.Range('E10').Value = '  -6.61%  '

# Row 11
$ws.Range('D11').Value = "'5.64"
$ws.Range('E11').Value = '  -5.06%  '

# Row 12
$ws.Range('D12').Value = "'0.442"
$ws.Range('E12').Value = '  -2.95%  '

# Row 13
$ws.Range('D13').Value = "'0.0000228"
$ws.Range('E13').Value = '  -5.40%  '

# Row 14
$ws.Range('D14').Value = "'34.57"
$ws.Range('E14').Value = '  -7.40%  '

# Row 15
$ws.Range('E15').Value = '  +1.47%  '

# Row 16
$ws.Range('D16').Value = "'3.495.40"
$ws.Range('E16').Value = '  -3.83%  '

# Row 17
$ws.Range('D17').Value = "'7.07"
$ws.Range('E17').Value = '  -2.86%  '

# Row 18
$ws.Range('D18').Value = "'62.334.13"
$ws.Range('E18').Value = '  -2.87%  '

# Row 19
$ws.Range('D19').Value = "'2.997.30"
$ws.Range('E19').Value = '  -3.85%  '

# Row 20
$ws.Range('D20').Value = "'456.10"
$ws.Range('E20').Value = '  -5.45%  '

# Row 21
$ws.Range('D21').Value = "'13.87"
$ws.Range('E21').Value = '  -4.87%  '

# Row 22
$ws.Range('D22').Value = "'0.679"
$ws.Range('E22').Value = '  -5.38%  '

# Row 23
$ws.Range('D23').Value = "'7.29"
$ws.Range('E23').Value = '  -4.58%  '

# Row 24
$ws.Range('D24').Value = "'79.98"
$ws.Range('E24').Value = '  -1.83%  '

# Row 25
$ws.Range('D25').Value = "'2.28"
$ws.Range('E25').Value = '  -7.90%  '

# Row 26
$ws.Range('D26').Value = "'12.27"
$ws.Range('E26').Value = '  -5.55%  '

# Row 27
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = "'10.01"
$ws.Range('E27').Value = '  -4.28%  '

# Row 28
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = "'1.00"
$ws.Range('E28').Value = '  -0.12%  '

# Row 29
$ws.Range('E29').Value = '  +0.05%  '

# Row 30
$ws.Range('D30').Value = "'7.18"
$ws.Range('E30').Value = '  -4.70%  '

# Row 31
$ws.Range('D31').Value = "'2.61"
$ws.Range('E31').Value = '  -3.44%  '

# Row 32
$ws.Range('D32').Value = "'2.09"
$ws.Range('E32').Value = '  -5.79%  '

# Row 33
$ws.Range('D33').Value = "'26.94"
$ws.Range('E33').Value = '  -2.07%  '

# Row 34
$ws.Range('D34').Value = "'0.107"
$ws.Range('E34').Value = '  -5.51%  '

# Row 35
$ws.Range('E35').Value = '  -3.84%  '

# Row 36
$ws.Range('D36').Value = "'0.0₃0785"
$ws.Range('E36').Value = '  -7.09%  '

# Row 37
$ws.Range('D37').Value = "'5.73"
$ws.Range('E37').Value = '  -5.51%  '

# Row 38
$ws.Range('D38').Value = "'2.11"
$ws.Range('E38').Value = '  -6.32%  '

# Row 39
$ws.Range('D39').Value = "'50.07"
$ws.Range('E39').Value = '  -1.99%  '

# Row 40
$ws.Range('D40').Value = "'9.04"
$ws.Range('E40').Value = '  -1.76%  '

# Row 41
$ws.Range('E41').Value = '  -11.44%  '

# Row 42
$ws.Range('D42').Value = "'410.54"
$ws.Range('E42').Value = '  -9.07%  '

# Row 43
$ws.Range('D43').Value = "'0.276"
$ws.Range('E43').Value = '  -5.72%  '

# Row 44
$ws.Range('D44').Value = "'0.111"
$ws.Range('E44').Value = '  -1.45%  '

# Row 45
$ws.Range('D45').Value = "'2.771.35"
$ws.Range('E45').Value = '  -2.76%  '

# Row 46
$ws.Range('D46').Value = "'0.0351"
$ws.Range('E46').Value = '  -4.32%  '

# Row 47
$ws.Range('D47').Value = "'38.17"
$ws.Range('E47').Value = '  -5.19%  '

# Row 48
$ws.Range('D48').Value = "'128.09"
$ws.Range('E48').Value = '  -1.73%  '

# Row 49
$ws.Range('E49').Value = '  +0.08%  '

# Row 50
$ws.Range('E50').Value = '  -2.37%  '

# Row 51
$ws.Range('D51').Value = "'23.73"
$ws.Range('E51').Value = '  -8.35%  '
